$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update value of C4 (code_class for TC001) and apply a scientific-notation
# number format (0.00E+00), matching the new style added to cellXfs.
$ws.Range("C4").Value = "5360e6355a"
$ws.Range("C4").NumberFormat = "0.00E+00"

# Reflect the new active selection on the sheet (was F7, now C4).
$ws.Range("C4").Select()
